# Commit: "changing to 5 years simulation"
# Update the "Coupling Parameters" sheet: change End Year from 2029 to 2025
# (Start Year 2020 + 5 years), and move the active cell selection to I7.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Coupling Parameters")
$ws.Activate()

# Set End Year (row 3, col B) to 2025 for a 5 year simulation
$ws.Range("B3").Value = 2025

# Update the active cell selection to I7
$ws.Range("I7").Select()
